# ng_db.xlsx edits:
#  - Row 3: liter "А4" -> "А3"
#  - Row 4: fill in the abk1 (2-этажное Здание АБК) record details
#  - Row 5: fill in the abk4 (3-этажное Здание АБК) record details
#  - Row heights for rows 4 and 5
#  - Active selection moves to I5
#  - Workbook no longer uses R1C1 reference style

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Helper: write a value as literal text (shared string), not auto-coerced to a
# number, even when the text looks numeric (e.g. "589.18"). We do this by
# entering it as a formula returning the literal string, then collapsing the
# formula down to its static value via copy / paste-special values. This
# preserves the cell's existing style (unlike flipping NumberFormat to "@").
function Set-TextValue {
    param($range, [string]$text)
    $escaped = $text -replace '"', '""'
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163) # xlPasteValues
}

# Row 3
Set-TextValue $ws.Range("F3") "А3"

# Row 4
$ws.Range("C4").Value = "2-этажное Здание АБК (1-очередь)  пл. 589,18 м2, 000000032"
Set-TextValue $ws.Range("H4") "589.18"
$ws.Range("A4:L4").RowHeight = 30

# Row 5
$ws.Range("C5").Value = "3-этажное Здание АБК  (2-я очередь), Ново-Гайвинская, 92 566,1 м2 инв.№117, 00-000007"
$ws.Range("D5").Value = 117
Set-TextValue $ws.Range("F5") "А2"
Set-TextValue $ws.Range("H5") "566.1"
$ws.Range("A5:L5").RowHeight = 45

# Selection moves to I5
$ws.Range("I5").Select()

# calcPr refMode="R1C1" is dropped -> switch the application back to A1 style
$excel.ReferenceStyle = 1
